# Initial support for cite-seq data
# Add two new columns (I: CITE-seq Library Index, J: CITE-seq Library Conc)
# to the ImportReadsetTemplate header row, matching the existing header style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new header values (this also extends dimension/spans automatically)
$ws.Range("I1").Value = "CITE-seq Library Index"
$ws.Range("J1").Value = "CITE-seq Library Conc"

# Copy the formatting (style incl. border) from an existing header cell (H1)
# onto A1 (which previously lacked the border used by the rest of the header
# row) as well as onto the two newly added header cells I1/J1.
$ws.Range("H1").Copy()
$ws.Range("A1:J1").PasteSpecial(-4122)

# Grow the header row to fit the now-larger (wrapped) header text.
$ws.Rows.Item(1).RowHeight = 47.25

# Reset the selection back to the default cell.
$ws.Range("A1").Select()
